$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "ID"
$ws.Range("F1").Value = "Amount"

$ws.Range("D2").Value = 9052046524
$ws.Range("E2").Value = 123456
$ws.Range("F2").Value = 12.34

$ws.Range("D3").Value = 9076565434
$ws.Range("E3").Value = 56456464
$ws.Range("F3").Value = 25

$ws.Columns("D:D").AutoFit()

$ws.Range("H4").Select()
